# Updates the "cryptos" price/volume table to the latest scrape.
#
# Column D (Price) and column E (Volume(1h)) are stored as plain text in
# this sheet (not numbers), because values like "58.630.50" use dots as
# thousands separators and the percentages carry literal "%" + padding
# spaces. Several of the new Price values (e.g. "529.67", "1.00") *do*
# look like ordinary decimal numbers though, and a bare
# `.Value = "529.67"` assignment lets Excel's input-parsing turn them into
# real numbers (losing the trailing zero / changing the cell type). For
# those cells we briefly force a Text number format, write the value, and
# then restore the cell's original ("Normal") style so no formatting
# artifact is left behind - matching the source XML, which keeps every
# D/E cell as a plain, unstyled inline string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "58.636.62"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "3.165.24"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "529.67"
$ws.Range("E5").Value = "  -0.14%  "
Set-TextValue "D6" "139.81"
$ws.Range("E6").Value = "  +1.23%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "0.539"
$ws.Range("E8").Value = "  +14.96%  "
$ws.Range("E9").Value = "  +0.04%  "
Set-TextValue "D10" "0.438"
$ws.Range("E10").Value = "  +5.97%  "
$ws.Range("E11").Value = "  +4.12%  "
$ws.Range("E12").Value = "  +2.58%  "
$ws.Range("D13").Value = "3.709.89"
$ws.Range("E13").Value = "  +2.00%  "
Set-TextValue "D14" "25.74"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").Value = "58.681.38"
$ws.Range("E16").Value = "  +1.75%  "
Set-TextValue "D17" "6.26"
$ws.Range("E17").Value = "  +3.72%  "
$ws.Range("D18").Value = "3.162.20"
$ws.Range("E18").Value = "  +2.31%  "
Set-TextValue "D19" "12.98"
$ws.Range("E19").Value = "  +2.41%  "
Set-TextValue "D20" "377.05"
$ws.Range("E20").Value = "  +4.57%  "
Set-TextValue "D21" "8.11"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("E22").Value = "  -0.14%  "
Set-TextValue "D23" "0.531"
$ws.Range("E23").Value = "  +5.34%  "
Set-TextValue "D24" "69.71"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  +0.24%  "
Set-TextValue "D27" "8.24"
$ws.Range("E27").Value = "  +12.64%  "
$ws.Range("D28").Value = "0.0₃0867"
$ws.Range("E28").Value = "  -0.02%  "
Set-TextValue "D29" "22.34"
$ws.Range("E29").Value = "  +4.67%  "
$ws.Range("E30").Value = "  +1.01%  "
Set-TextValue "D31" "6.05"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +4.34%  "
Set-TextValue "D35" "156.93"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  +4.77%  "
Set-TextValue "D37" "25.12"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "2.680.10"
$ws.Range("E38").Value = "  +7.68%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D39" "0.0696"
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "1.69"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("E41").Value = "  +7.22%  "
Set-TextValue "D42" "0.723"
$ws.Range("E42").Value = "  +4.06%  "
Set-TextValue "D43" "39.14"
$ws.Range("E43").Value = "  +4.14%  "
$ws.Range("E44").Value = "  +7.39%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "3.206.62"
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("E47").Value = "  +14.40%  "
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  -0.27%  "
Set-TextValue "D50" "20.12"
$ws.Range("E50").Value = "  +1.86%  "
Set-TextValue "D51" "0.751"
$ws.Range("E51").Value = "  +1.70%  "
